# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Wed Oct  4 22:38:50 UTC 2023 with GitHub Actions".
#
# Two row pairs were also re-sorted by the refresh (their rank order
# changed), which shows up as the B/C/D/E cell contents of those rows
# being swapped between the two row numbers:
#   rows 19/20  -> ShibaInu/Chainlink swap places
#   rows 42/43  -> PaxDollar/FraxShare swap places
#   rows 44/45  -> mCoin/Aave swap places

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) that hold values Excel would otherwise auto-convert
# into numbers (stripping trailing zeros / thousands-style dots, switching to
# scientific notation, etc.) are first forced to Text format so the literal
# string from the source feed is preserved exactly, just like the original
# inline string cells. (NumberFormat is set per-cell: union ranges only
# apply to their first area.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "27.735.54"
$ws.Range("E2").Value = "  +1.48%  "

# Row 3
$ws.Range("D3").Value = "1.646.18"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "213.59"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").Value = "0.533"
$ws.Range("E6").Value = "  +4.10%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "23.26"
$ws.Range("E8").Value = "  -1.06%  "

# Row 9
$ws.Range("E9").Value = "  +0.57%  "

# Row 10
$ws.Range("E10").Value = "  +0.13%  "

# Row 11
$ws.Range("E11").Value = "  +0.27%  "

# Row 12
$ws.Range("D12").Value = "1.877.76"

# Row 13
$ws.Range("D13").Value = "1.646.98"
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$ws.Range("E14").Value = "  -0.72%  "

# Row 15
$ws.Range("D15").Value = "0.565"
$ws.Range("E15").Value = "  -0.76%  "

# Row 16
$ws.Range("D16").Value = "64.32"
$ws.Range("E16").Value = "  -1.74%  "

# Row 17
$ws.Range("D17").Value = "27.721.75"
$ws.Range("E17").Value = "  +1.40%  "

# Row 18
$ws.Range("D18").Value = "232.10"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +3.80%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0726"
$ws.Range("E20").Value = "  +0.10%  "

# Row 21
$ws.Range("E21").Value = "  -0.02%  "

# Row 22
$ws.Range("D22").Value = "4.33"
$ws.Range("E22").Value = "  -0.92%  "

# Row 23
$ws.Range("D23").Value = "10.12"
$ws.Range("E23").Value = "  +7.04%  "

# Row 24
$ws.Range("E24").Value = "  -3.84%  "

# Row 25
$ws.Range("D25").Value = "150.06"

# Row 26
$ws.Range("D26").Value = "6.99"
$ws.Range("E26").Value = "  -1.71%  "

# Row 27
$ws.Range("D27").Value = "0.113"
$ws.Range("E27").Value = "  +0.86%  "

# Row 28
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("D29").Value = "15.67"

# Row 30
$ws.Range("E30").Value = "  -0.03%  "

# Row 31
$ws.Range("D31").Value = "0.0488"
$ws.Range("E31").Value = "  -1.75%  "

# Row 32
$ws.Range("E32").Value = "  +0.58%  "

# Row 33
$ws.Range("D33").Value = "1.446.11"
$ws.Range("E33").Value = "  +1.03%  "

# Row 34
$ws.Range("D34").Value = "3.17"
$ws.Range("E34").Value = "  +0.83%  "

# Row 35
$ws.Range("E35").Value = "  +1.77%  "

# Row 36
$ws.Range("D36").Value = "2.34"
$ws.Range("E36").Value = "  -1.03%  "

# Row 37
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("D38").Value = "0.886"
$ws.Range("E38").Value = "  -2.25%  "

# Row 39
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("E40").Value = "  +14.31%  "

# Row 41
$ws.Range("E41").Value = "  -1.83%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.71"
$ws.Range("E42").Value = "  +2.49%  "

# Row 43
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "66.35"
$ws.Range("E44").Value = "  +1.98%  "

# Row 45
$ws.Range("B45").Value = "mCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  -0.70%  "

# Row 46
$ws.Range("E46").Value = "  +1.71%  "

# Row 47
$ws.Range("D47").Value = "1.786.99"
$ws.Range("E47").Value = "  -0.36%  "

# Row 48
$ws.Range("D48").Value = "1.72"
$ws.Range("E48").Value = "  +3.29%  "

# Row 49
$ws.Range("D49").Value = "86.53"
$ws.Range("E49").Value = "  -1.61%  "

# Row 50
$ws.Range("E50").Value = "  +1.65%  "

# Row 51
$ws.Range("D51").Value = "0.0995"
$ws.Range("E51").Value = "  -1.78%  "
